$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.16520000000001
$ws.Range("E4").Value = 16.17839999999999
$ws.Range("D7").Value = -7.293200000000001
$ws.Range("B8").Value = 5.838299999999997
$ws.Range("B10").Value = 4.7911
$ws.Range("E11").Value = 15.95330000000001
$ws.Range("B12").Value = 5.152300000000002
$ws.Range("D14").Value = -8.439400000000004
$ws.Range("E14").Value = 16.34709999999999
$ws.Range("D15").Value = -7.845100000000002
$ws.Range("B18").Value = 5.628599999999995
$ws.Range("D18").Value = -8.878299999999989
$ws.Range("E18").Value = 16.48789999999998
$ws.Range("E19").Value = 16.55550000000001
$ws.Range("D20").Value = -7.995099999999995
$ws.Range("E21").Value = 16.95470000000001
$ws.Range("B25").Value = 5.890000000000001
$ws.Range("E27").Value = 16.54599999999999
$ws.Range("D29").Value = -6.8985
$ws.Range("D30").Value = -7.924399999999995
$ws.Range("D31").Value = -7.555499999999999
$ws.Range("E31").Value = 16.8596
$ws.Range("D35").Value = -8.562099999999994
$ws.Range("B37").Value = 8.5387
$ws.Range("E38").Value = 16.48699999999999
$ws.Range("D40").Value = -8.525599999999994
$ws.Range("E42").Value = 16.36109999999999
$ws.Range("D44").Value = -7.540300000000001
$ws.Range("E44").Value = 16.7372
$ws.Range("E47").Value = 16.56819999999999
$ws.Range("D50").Value = -7.983599999999998
$ws.Range("D54").Value = -7.932200000000005
$ws.Range("B55").Value = 5.518699999999995
$ws.Range("E56").Value = 16.2576
$ws.Range("E58").Value = 15.91570000000002
$ws.Range("E65").Value = 17.19480000000001
$ws.Range("B68").Value = 5.095899999999999
$ws.Range("D68").Value = -7.057599999999993
$ws.Range("E73").Value = 17.40220000000001
$ws.Range("D76").Value = -7.369199999999998
$ws.Range("B77").Value = 9.680500000000004
$ws.Range("B78").Value = 9.656700000000003
$ws.Range("B79").Value = 9.226400000000003
$ws.Range("B80").Value = 9.025499999999999
$ws.Range("B81").Value = 5.038700000000005
$ws.Range("B82").Value = 6.268800000000002
$ws.Range("B84").Value = 5.449300000000001
$ws.Range("D87").Value = -8.011099999999995
$ws.Range("D88").Value = -7.621599999999998
$ws.Range("E90").Value = 16.5427
$ws.Range("D92").Value = -7.330200000000002
$ws.Range("E92").Value = 16.81800000000001
$ws.Range("E94").Value = 19.25020000000002
$ws.Range("E95").Value = 18.19730000000002
$ws.Range("D96").Value = -7.7912
$ws.Range("D98").Value = -8.487200000000009
$ws.Range("B101").Value = 8.7323
$ws.Range("D101").Value = -7.697100000000002
$ws.Range("E101").Value = 16.4863
$ws.Range("B102").Value = 8.6966
$ws.Range("D102").Value = -7.559700000000001
